# Applies the "feat: add 2022-Q4 data" edit:
#   1. Insert a new worksheet "2022-Q4" right after "总计", holding the
#      new quarter's fund-holdings table (shifts every later Q-sheet
#      one tab to the right).
#   2. Prepend a "2022-Q4" row to the "总计" summary sheet, and append a
#      trailing "2020-Q4" row (the summary table grows from 8 to 9 data rows).

function Set-TextValue($cell, [string]$val) {
    # Force literal-text storage so numeric-looking strings ("40.95",
    # "217024", "010430" ...) keep their original formatting / leading
    # zeros instead of being coerced into numbers by Excel's smart entry.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q4" worksheet, positioned right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# Match the page-margin defaults used by every other sheet in the workbook.
$q4.PageSetup.LeftMargin = 0.75 * 72
$q4.PageSetup.RightMargin = 0.75 * 72
$q4.PageSetup.TopMargin = 1 * 72
$q4.PageSetup.BottomMargin = 1 * 72
$q4.PageSetup.HeaderMargin = 0.5 * 72
$q4.PageSetup.FooterMargin = 0.5 * 72

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4.Cells.Item(1, $col).Value = $headers[$col - 2]
}
# Match the header-row styling used on every other quarter sheet.
$totalSheet.Range("B1").Copy() | Out-Null
$q4.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$rows = @(
    @("217024", "招商安盈债券A",                   "40.95", "20.32", "1.93", "0.7903", 3),
    @("519198", "万家颐和灵活配置混合A",            "9.08",  "93.26", "5.70", "0.5176", 8),
    @("014887", "招商安福1年定期开放债券",           "17.78", "33.59", "1.49", "0.2649", 6),
    @("008979", "万家民丰回报一年持有期混合",         "16.63", "29.23", "1.19", "0.1979", 10),
    @("010430", "招商安阳债券A",                    "16.90", "20.35", "1.09", "0.1842", 6),
    @("519183", "万家双引擎灵活配置混合",             "1.97",  "93.92", "7.75", "0.1527", 3),
    @("016513", "招商安嘉债券",                     "16.27", "20.17", "0.76", "0.1237", 8),
    @("016620", "万家颐和灵活配置混合C",             "1.23",  "93.26", "5.70", "0.0701", 8),
    @("011018", "景顺长城安泽回报一年持有期混合A",     "3.00",  "37.02", "1.30", "0.0390", 10),
    @("519197", "万家颐达灵活配置混合",               "1.24",  "44.22", "2.55", "0.0316", 8),
    @("014768", "景顺华城稳健6月持有混合C",           "1.61",  "22.61", "1.10", "0.0177", 5),
    @("012977", "瑞达鑫红量化6个月持有混合A",         "0.35",  "94.66", "4.91", "0.0172", 9),
    @("014767", "景顺华城稳健6月持有混合A",           "1.10",  "22.61", "1.10", "0.0121", 5),
    @("012978", "瑞达鑫红量化6个月持有混合C",         "0.09",  "94.66", "4.91", "0.0044", 9),
    @("011019", "景顺长城安泽回报一年持有期混合C",     "0.15",  "37.02", "1.30", "0.0020", 10),
    @("010431", "招商安阳债券C",                    "0.11",  "20.35", "1.09", "0.0012", 6),
    @("012233", "招商安盈债券C",                    "0.01",  "20.32", "1.93", "0.0002", 3)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $r - 2
    Set-TextValue $q4.Cells.Item($r, 2) $row[0]
    Set-TextValue $q4.Cells.Item($r, 3) $row[1]
    Set-TextValue $q4.Cells.Item($r, 4) $row[2]
    Set-TextValue $q4.Cells.Item($r, 5) $row[3]
    Set-TextValue $q4.Cells.Item($r, 6) $row[4]
    Set-TextValue $q4.Cells.Item($r, 7) $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Match column-A index styling used on every other quarter sheet.
$totalSheet.Range("A2").Copy() | Out-Null
$q4.Range("A2:A18").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: new 2022-Q4 row up top, and the
#    table grows a trailing 2020-Q4 row (data rows 2-9).
# ---------------------------------------------------------------------
$summary = @(
    @("2022-Q4", 17, 2.43),
    @("2022-Q3", 63, 22.68),
    @("2022-Q2", 72, 19.36),
    @("2022-Q1", 53, 8.99),
    @("2021-Q4", 8, 1.26),
    @("2021-Q3", 17, 3.53),
    @("2021-Q2", 1, 2.06),
    @("2020-Q4", 1, 0.06)
)

$r = 2
foreach ($d in $summary) {
    $totalSheet.Cells.Item($r, 2).Value = $d[0]
    $totalSheet.Cells.Item($r, 3).Value = $d[1]
    $totalSheet.Cells.Item($r, 4).Value = $d[2]
    $r++
}

# New row 9 needs the same column-A index styling as rows 2-8.
$totalSheet.Cells.Item(9, 1).Value = 7
$totalSheet.Range("A2").Copy() | Out-Null
$totalSheet.Range("A9").PasteSpecial(-4122) | Out-Null

# Restore "总计" as the active/selected sheet (it was selected before this
# edit; adding/populating "2022-Q4" would otherwise leave it focused).
$totalSheet.Activate()
$totalSheet.Range("A1").Select() | Out-Null
